$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# C1 / D1: trim the boxed border so the merged-header underline (B1:D1)
# doesn't double up on the inner cell edges. C1 keeps top+bottom only;
# D1 (right-most cell of the merge) keeps top+right+bottom.
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.LineStyle = 1
$c1.Borders(7).LineStyle = -4142
$c1.Borders(10).LineStyle = -4142

$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.LineStyle = 1
$d1.Borders(7).LineStyle = -4142

# Rename header column from "fedcore" to "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

$c1b = $ws2.Range("C1")
$c1b.ClearFormats()
$c1b.Borders.LineStyle = 1
$c1b.Borders(7).LineStyle = -4142
$c1b.Borders(10).LineStyle = -4142

$d1b = $ws2.Range("D1")
$d1b.ClearFormats()
$d1b.Borders.LineStyle = 1
$d1b.Borders(7).LineStyle = -4142

$f1b = $ws2.Range("F1")
$f1b.ClearFormats()
$f1b.Borders.LineStyle = 1
$f1b.Borders(7).LineStyle = -4142
$f1b.Borders(10).LineStyle = -4142

$g1b = $ws2.Range("G1")
$g1b.ClearFormats()
$g1b.Borders.LineStyle = 1
$g1b.Borders(7).LineStyle = -4142

# Rename header columns from "fedcore" to "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell at G5
$ws2.Range("G5").ClearContents()
